$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old rows 11 and 12 (being removed)
$ws.Range("A11:B12").ClearContents()

# Update header
$ws.Range("A1").Value = "Cluster name"
$ws.Range("B1").Value = "Active cases"

# Update data rows A2:B10
$data = @(
    @("3155 Westmont Aged Care Services Baranduda", 10),
    @("3642 Fronditha Care Aged Care Clayton South", 35),
    @("4314 Estia Health Altona Meadows", 10),
    @("Confirmed Omicron Sircuit Bar Fitzroy", 20),
    @("Confirmed Omicron Variant The Peel Hotel Collingwood", 23),
    @("Diamond Valley Pork and Baxters Pork Laverton North", 36),
    @("Mercure Welcome Melbourne", 13),
    @("Pullman Melbourne on Swanston Melbourne", 10),
    @("Werribee Mercy Hospital Emergency Department", 17)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
